$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.965
$ws.Range("B9").Value = 5.552
$ws.Range("C12").Value = -11.601
$ws.Range("C14").Value = -12.845
$ws.Range("B18").Value = 5.103
$ws.Range("B20").Value = 6.85
$ws.Range("C26").Value = -13.131
$ws.Range("B27").Value = 5.526000000000001
$ws.Range("C27").Value = -13.664
$ws.Range("C29").Value = -12.233
$ws.Range("B35").Value = 7.657999999999999
$ws.Range("C37").Value = -13.426
$ws.Range("C38").Value = -13.313
$ws.Range("C51").Value = -11.591
$ws.Range("C52").Value = -11.621
$ws.Range("C55").Value = -13.752
$ws.Range("B69").Value = 5.667
$ws.Range("C69").Value = -10.732
$ws.Range("C70").Value = -11.645
$ws.Range("B76").Value = 6.723999999999999
$ws.Range("B78").Value = 7.230999999999999
$ws.Range("C81").Value = -13.418
$ws.Range("B82").Value = 5.061
$ws.Range("B83").Value = 5.404999999999999
$ws.Range("C83").Value = -13.668
$ws.Range("B93").Value = 5.628
$ws.Range("C102").Value = -13.419
